# Chap 8 Royal Company class case - fill in the Master Budget input cells
# and touch a couple of cosmetic items (number format on the MOH line,
# and the last-saved selection) to match the target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget w.o Data")

# ---- Budget assumption inputs (rows 19-43, column B) ----------------------
$ws.Range("B19").Value = 0.5
$ws.Range("B20").Value = 0.5
$ws.Range("B21").Value = 12000
$ws.Range("B22").Value = 0.05
$ws.Range("B23").Value = 10
$ws.Range("B24").Value = 40
$ws.Range("B26").Value = 20
$ws.Range("B27").Value = 30000
$ws.Range("B28").Value = 20000
# Manufacturing overhead rate per hour is derived, not typed in directly
$ws.Range("B29").Formula = "=E105/SUM(B101:D101)"
$ws.Range("B30").Value = 0.5
$ws.Range("B31").Value = 60000
$ws.Range("B32").Value = 10000
$ws.Range("B33").Value = 30000
$ws.Range("B34").Value = 49000
$ws.Range("B35").Value = 40000
$ws.Range("B36").Value = 143700
$ws.Range("B37").Value = 48300
$ws.Range("B38").Value = 75000
$ws.Range("B39").Value = 0.16
$ws.Range("B40").Value = 50000
$ws.Range("B41").Value = 150000
$ws.Range("B42").Value = 248650
$ws.Range("B43").Value = 175000

# ---- Cosmetic: direct labor cost line (row 91) drops the $ sign -----------
# (keeps the existing font / right-center alignment, just swaps the number
# format from a currency accounting format to a plain accounting format)
$ws.Range("B91:E91").NumberFormat = "_-* #,##0_-;\-* #,##0_-;_-* ""-""??_-;_-@_-"

# ---- Recalculate so every formula's cached <v> is fresh --------------------
$excel.Calculate()

# ---- Restore the sheet's last-saved cursor position ------------------------
$ws.Activate()
$ws.Range("D166").Select()
